$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 16.71895933333333
$ws.Cells.Item(2, 8).Value = 50.156878
$ws.Cells.Item(2, 9).Value = 0.02912144738161902
$ws.Cells.Item(2, 10).Value = 0.03059269312988411
$ws.Cells.Item(2, 13).Value = 0.8151449999999999
$ws.Cells.Item(2, 14).Value = 2.445435
$ws.Cells.Item(2, 15).Value = 0.1271069095499719
$ws.Cells.Item(2, 16).Value = 0.1371035811308388
$ws.Cells.Item(2, 17).Value = 13.62837610577
$ws.Cells.Item(2, 18).Value = 122.65538495193
$ws.Cells.Item(2, 19).Value = 0.003701537178299713
$ws.Cells.Item(2, 20).Value = 0.004194367784543922
# Row 3
$ws.Cells.Item(3, 7).Value = 16.71895933333333
$ws.Cells.Item(3, 8).Value = 50.156878
$ws.Cells.Item(3, 9).Value = 0.02912144738161902
$ws.Cells.Item(3, 10).Value = 0.03059269312988411
$ws.Cells.Item(3, 15).Value = 0.4802730342501803
$ws.Cells.Item(3, 16).Value = 0.5180454245123947
$ws.Cells.Item(3, 17).Value = 51.49477370974489
$ws.Cells.Item(3, 18).Value = 463.452963387704
$ws.Cells.Item(3, 19).Value = 0.01398624589572713
$ws.Cells.Item(3, 20).Value = 0.01584840469944824
# Row 4
$ws.Cells.Item(4, 7).Value = 16.71895933333333
$ws.Cells.Item(4, 8).Value = 50.156878
$ws.Cells.Item(4, 9).Value = 0.02912144738161902
$ws.Cells.Item(4, 10).Value = 0.03059269312988411
$ws.Cells.Item(4, 13).Value = 0.5185940000000001
$ws.Cells.Item(4, 14).Value = 1.555782
$ws.Cells.Item(4, 15).Value = 0.08086522109705406
$ws.Cells.Item(4, 16).Value = 0.08722508823947427
$ws.Cells.Item(4, 17).Value = 8.67035199651067
$ws.Cells.Item(4, 18).Value = 78.033167968596
$ws.Cells.Item(4, 19).Value = 0.002354912281180848
$ws.Cells.Item(4, 20).Value = 0.0026684503577373
# Row 5
$ws.Cells.Item(5, 7).Value = 16.71895933333333
$ws.Cells.Item(5, 8).Value = 50.156878
$ws.Cells.Item(5, 9).Value = 0.02912144738161902
$ws.Cells.Item(5, 10).Value = 0.03059269312988411
$ws.Cells.Item(5, 13).Value = 1.402793
$ws.Cells.Item(5, 14).Value = 2.805586
$ws.Cells.Item(5, 15).Value = 0.2187398352051889
$ws.Cells.Item(5, 16).Value = 0.1572954863942594
$ws.Cells.Item(5, 17).Value = 23.45323912008467
$ws.Cells.Item(5, 18).Value = 140.719434720508
$ws.Cells.Item(5, 19).Value = 0.006370020601191924
$ws.Cells.Item(5, 20).Value = 0.004812092545975438
# Row 6
$ws.Cells.Item(6, 7).Value = 16.71895933333333
$ws.Cells.Item(6, 8).Value = 50.156878
$ws.Cells.Item(6, 9).Value = 0.02912144738161902
$ws.Cells.Item(6, 10).Value = 0.03059269312988411
$ws.Cells.Item(6, 13).Value = 0.5965113333333333
$ws.Cells.Item(6, 14).Value = 1.789534
$ws.Cells.Item(6, 15).Value = 0.09301499989760488
$ws.Cells.Item(6, 16).Value = 0.1003304197230327
$ws.Cells.Item(6, 17).Value = 9.973048723872443
$ws.Cells.Item(6, 18).Value = 89.757438514852
$ws.Cells.Item(6, 19).Value = 0.002708731425219399
$ws.Cells.Item(6, 20).Value = 0.003069377742179213
# Row 7
$ws.Cells.Item(7, 9).Value = 0.2708539632042961
$ws.Cells.Item(7, 10).Value = 0.2845377865576845
$ws.Cells.Item(7, 13).Value = 0.8151449999999999
$ws.Cells.Item(7, 14).Value = 2.445435
$ws.Cells.Item(7, 15).Value = 0.1271069095499719
$ws.Cells.Item(7, 16).Value = 0.1371035811308388
$ws.Cells.Item(7, 17).Value = 126.755364591405
$ws.Cells.Item(7, 18).Value = 1140.798281322645
$ws.Cells.Item(7, 19).Value = 0.03442741020225987
$ws.Cells.Item(7, 20).Value = 0.0390111495041008
# Row 8
$ws.Cells.Item(8, 9).Value = 0.2708539632042961
$ws.Cells.Item(8, 10).Value = 0.2845377865576845
$ws.Cells.Item(8, 15).Value = 0.4802730342501803
$ws.Cells.Item(8, 16).Value = 0.5180454245123947
$ws.Cells.Item(8, 19).Value = 0.1300838547468139
$ws.Cells.Item(8, 20).Value = 0.1474034984270928
# Row 9
$ws.Cells.Item(9, 9).Value = 0.2708539632042961
$ws.Cells.Item(9, 10).Value = 0.2845377865576845
$ws.Cells.Item(9, 13).Value = 0.5185940000000001
$ws.Cells.Item(9, 14).Value = 1.555782
$ws.Cells.Item(9, 15).Value = 0.08086522109705406
$ws.Cells.Item(9, 16).Value = 0.08722508823947427
$ws.Cells.Item(9, 17).Value = 80.64156873306602
$ws.Cells.Item(9, 18).Value = 725.7741185975941
$ws.Cells.Item(9, 19).Value = 0.02190266561952875
$ws.Cells.Item(9, 20).Value = 0.02481883353995873
# Row 10
$ws.Cells.Item(10, 9).Value = 0.2708539632042961
$ws.Cells.Item(10, 10).Value = 0.2845377865576845
$ws.Cells.Item(10, 13).Value = 1.402793
$ws.Cells.Item(10, 14).Value = 2.805586
$ws.Cells.Item(10, 15).Value = 0.2187398352051889
$ws.Cells.Item(10, 16).Value = 0.1572954863942594
$ws.Cells.Item(10, 17).Value = 218.134857186477
$ws.Cells.Item(10, 18).Value = 1308.809143118862
$ws.Cells.Item(10, 19).Value = 0.05924655127598002
$ws.Cells.Item(10, 20).Value = 0.04475650953413694
# Row 11
$ws.Cells.Item(11, 9).Value = 0.2708539632042961
$ws.Cells.Item(11, 10).Value = 0.2845377865576845
$ws.Cells.Item(11, 13).Value = 0.5965113333333333
$ws.Cells.Item(11, 14).Value = 1.789534
$ws.Cells.Item(11, 15).Value = 0.09301499989760488
$ws.Cells.Item(11, 16).Value = 0.1003304197230327
$ws.Cells.Item(11, 17).Value = 92.757744376242
$ws.Cells.Item(11, 18).Value = 834.8196993861779
$ws.Cells.Item(11, 19).Value = 0.02519348135971347
$ws.Cells.Item(11, 20).Value = 0.02854779555239519
# Row 12
$ws.Cells.Item(12, 7).Value = 194.8548433333333
$ws.Cells.Item(12, 8).Value = 584.56453
$ws.Cells.Item(12, 9).Value = 0.3394024086099587
$ws.Cells.Item(12, 10).Value = 0.3565493705749576
$ws.Cells.Item(12, 13).Value = 0.8151449999999999
$ws.Cells.Item(12, 14).Value = 2.445435
$ws.Cells.Item(12, 15).Value = 0.1271069095499719
$ws.Cells.Item(12, 16).Value = 0.1371035811308388
$ws.Cells.Item(12, 17).Value = 158.83495126895
$ws.Cells.Item(12, 18).Value = 1429.51456142055
$ws.Cells.Item(12, 19).Value = 0.04314039125222862
$ws.Cells.Item(12, 20).Value = 0.04888419555577322
# Row 13
$ws.Cells.Item(13, 7).Value = 194.8548433333333
$ws.Cells.Item(13, 8).Value = 584.56453
$ws.Cells.Item(13, 9).Value = 0.3394024086099587
$ws.Cells.Item(13, 10).Value = 0.3565493705749576
$ws.Cells.Item(13, 15).Value = 0.4802730342501803
$ws.Cells.Item(13, 16).Value = 0.5180454245123947
$ws.Cells.Item(13, 17).Value = 600.1573341764489
$ws.Cells.Item(13, 18).Value = 5401.41600758804
$ws.Cells.Item(13, 19).Value = 0.1630058246149244
$ws.Cells.Item(13, 20).Value = 0.184708770039131
# Row 14
$ws.Cells.Item(14, 7).Value = 194.8548433333333
$ws.Cells.Item(14, 8).Value = 584.56453
$ws.Cells.Item(14, 9).Value = 0.3394024086099587
$ws.Cells.Item(14, 10).Value = 0.3565493705749576
$ws.Cells.Item(14, 13).Value = 0.5185940000000001
$ws.Cells.Item(14, 14).Value = 1.555782
$ws.Cells.Item(14, 15).Value = 0.08086522109705406
$ws.Cells.Item(14, 16).Value = 0.08722508823947427
$ws.Cells.Item(14, 17).Value = 101.0505526236067
$ws.Cells.Item(14, 18).Value = 909.4549736124601
$ws.Cells.Item(14, 19).Value = 0.027445850813117
$ws.Cells.Item(14, 20).Value = 0.03110005031012968
# Row 15
$ws.Cells.Item(15, 7).Value = 194.8548433333333
$ws.Cells.Item(15, 8).Value = 584.56453
$ws.Cells.Item(15, 9).Value = 0.3394024086099587
$ws.Cells.Item(15, 10).Value = 0.3565493705749576
$ws.Cells.Item(15, 13).Value = 1.402793
$ws.Cells.Item(15, 14).Value = 2.805586
$ws.Cells.Item(15, 15).Value = 0.2187398352051889
$ws.Cells.Item(15, 16).Value = 0.1572954863942594
$ws.Cells.Item(15, 17).Value = 273.3410102440967
$ws.Cells.Item(15, 18).Value = 1640.04606146458
$ws.Cells.Item(15, 19).Value = 0.07424082692758656
$ws.Cells.Item(15, 20).Value = 0.05608360666815498
# Row 16
$ws.Cells.Item(16, 7).Value = 194.8548433333333
$ws.Cells.Item(16, 8).Value = 584.56453
$ws.Cells.Item(16, 9).Value = 0.3394024086099587
$ws.Cells.Item(16, 10).Value = 0.3565493705749576
$ws.Cells.Item(16, 13).Value = 0.5965113333333333
$ws.Cells.Item(16, 14).Value = 1.789534
$ws.Cells.Item(16, 15).Value = 0.09301499989760488
$ws.Cells.Item(16, 16).Value = 0.1003304197230327
$ws.Cells.Item(16, 17).Value = 116.2331224032244
$ws.Cells.Item(16, 18).Value = 1046.09810162902
$ws.Cells.Item(16, 19).Value = 0.03156951500210217
$ws.Cells.Item(16, 20).Value = 0.03577274800176863
# Row 17
$ws.Cells.Item(17, 7).Value = 82.82950199999999
$ws.Cells.Item(17, 8).Value = 165.659004
$ws.Cells.Item(17, 9).Value = 0.1442742299952585
$ws.Cells.Item(17, 10).Value = 0.1010420758958371
$ws.Cells.Item(17, 13).Value = 0.8151449999999999
$ws.Cells.Item(17, 14).Value = 2.445435
$ws.Cells.Item(17, 15).Value = 0.1271069095499719
$ws.Cells.Item(17, 16).Value = 0.1371035811308388
$ws.Cells.Item(17, 17).Value = 67.51805440778999
$ws.Cells.Item(17, 18).Value = 405.1083264467399
$ws.Cells.Item(17, 19).Value = 0.01833825150239916
$ws.Cells.Item(17, 20).Value = 0.01385323045021328
# Row 18
$ws.Cells.Item(18, 7).Value = 82.82950199999999
$ws.Cells.Item(18, 8).Value = 165.659004
$ws.Cells.Item(18, 9).Value = 0.1442742299952585
$ws.Cells.Item(18, 10).Value = 0.1010420758958371
$ws.Cells.Item(18, 15).Value = 0.4802730342501803
$ws.Cells.Item(18, 16).Value = 0.5180454245123947
$ws.Cells.Item(18, 17).Value = 255.116743628712
$ws.Cells.Item(18, 18).Value = 1530.700461772272
$ws.Cells.Item(18, 19).Value = 0.06929102220393117
$ws.Cells.Item(18, 20).Value = 0.05234438510107255
# Row 19
$ws.Cells.Item(19, 7).Value = 82.82950199999999
$ws.Cells.Item(19, 8).Value = 165.659004
$ws.Cells.Item(19, 9).Value = 0.1442742299952585
$ws.Cells.Item(19, 10).Value = 0.1010420758958371
$ws.Cells.Item(19, 13).Value = 0.5185940000000001
$ws.Cells.Item(19, 14).Value = 1.555782
$ws.Cells.Item(19, 15).Value = 0.08086522109705406
$ws.Cells.Item(19, 16).Value = 0.08722508823947427
$ws.Cells.Item(19, 17).Value = 42.95488276018801
$ws.Cells.Item(19, 18).Value = 257.729296561128
$ws.Cells.Item(19, 19).Value = 0.01166676750717381
$ws.Cells.Item(19, 20).Value = 0.008813403985914049
# Row 20
$ws.Cells.Item(20, 7).Value = 82.82950199999999
$ws.Cells.Item(20, 8).Value = 165.659004
$ws.Cells.Item(20, 9).Value = 0.1442742299952585
$ws.Cells.Item(20, 10).Value = 0.1010420758958371
$ws.Cells.Item(20, 13).Value = 1.402793
$ws.Cells.Item(20, 14).Value = 2.805586
$ws.Cells.Item(20, 15).Value = 0.2187398352051889
$ws.Cells.Item(20, 16).Value = 0.1572954863942594
$ws.Cells.Item(20, 17).Value = 116.192645599086
$ws.Cells.Item(20, 18).Value = 464.7705823963439
$ws.Cells.Item(20, 19).Value = 0.03155852129351837
$ws.Cells.Item(20, 20).Value = 0.01589346247432137
# Row 21
$ws.Cells.Item(21, 7).Value = 82.82950199999999
$ws.Cells.Item(21, 8).Value = 165.659004
$ws.Cells.Item(21, 9).Value = 0.1442742299952585
$ws.Cells.Item(21, 10).Value = 0.1010420758958371
$ws.Cells.Item(21, 13).Value = 0.5965113333333333
$ws.Cells.Item(21, 14).Value = 1.789534
$ws.Cells.Item(21, 15).Value = 0.09301499989760488
$ws.Cells.Item(21, 16).Value = 0.1003304197230327
$ws.Cells.Item(21, 17).Value = 49.40873667735599
$ws.Cells.Item(21, 18).Value = 296.4524200641359
$ws.Cells.Item(21, 19).Value = 0.01341966748823599
$ws.Cells.Item(21, 20).Value = 0.01013759388431587
# Row 22
$ws.Cells.Item(22, 7).Value = 124.2078576666667
$ws.Cells.Item(22, 8).Value = 372.623573
$ws.Cells.Item(22, 9).Value = 0.2163479508088675
$ws.Cells.Item(22, 10).Value = 0.2272780738416368
$ws.Cells.Item(22, 13).Value = 0.8151449999999999
$ws.Cells.Item(22, 14).Value = 2.445435
$ws.Cells.Item(22, 15).Value = 0.1271069095499719
$ws.Cells.Item(22, 16).Value = 0.1371035811308388
$ws.Cells.Item(22, 17).Value = 101.247414137695
$ws.Cells.Item(22, 18).Value = 911.2267272392548
$ws.Cells.Item(22, 19).Value = 0.02749931941478449
$ws.Cells.Item(22, 20).Value = 0.03116063783620764
# Row 23
$ws.Cells.Item(23, 7).Value = 124.2078576666667
$ws.Cells.Item(23, 8).Value = 372.623573
$ws.Cells.Item(23, 9).Value = 0.2163479508088675
$ws.Cells.Item(23, 10).Value = 0.2272780738416368
$ws.Cells.Item(23, 15).Value = 0.4802730342501803
$ws.Cells.Item(23, 16).Value = 0.5180454245123947
$ws.Cells.Item(23, 17).Value = 382.5630169914404
$ws.Cells.Item(23, 18).Value = 3443.067152922964
$ws.Cells.Item(23, 19).Value = 0.1039060867887836
$ws.Cells.Item(23, 20).Value = 0.1177403662456501
# Row 24
$ws.Cells.Item(24, 7).Value = 124.2078576666667
$ws.Cells.Item(24, 8).Value = 372.623573
$ws.Cells.Item(24, 9).Value = 0.2163479508088675
$ws.Cells.Item(24, 10).Value = 0.2272780738416368
$ws.Cells.Item(24, 13).Value = 0.5185940000000001
$ws.Cells.Item(24, 14).Value = 1.555782
$ws.Cells.Item(24, 15).Value = 0.08086522109705406
$ws.Cells.Item(24, 16).Value = 0.08722508823947427
$ws.Cells.Item(24, 17).Value = 64.41344973878734
$ws.Cells.Item(24, 18).Value = 579.721047649086
$ws.Cells.Item(24, 19).Value = 0.01749502487605365
$ws.Cells.Item(24, 20).Value = 0.01982435004573452
# Row 25
$ws.Cells.Item(25, 7).Value = 124.2078576666667
$ws.Cells.Item(25, 8).Value = 372.623573
$ws.Cells.Item(25, 9).Value = 0.2163479508088675
$ws.Cells.Item(25, 10).Value = 0.2272780738416368
$ws.Cells.Item(25, 13).Value = 1.402793
$ws.Cells.Item(25, 14).Value = 2.805586
$ws.Cells.Item(25, 15).Value = 0.2187398352051889
$ws.Cells.Item(25, 16).Value = 0.1572954863942594
$ws.Cells.Item(25, 17).Value = 174.2379132797963
$ws.Cells.Item(25, 18).Value = 1045.427479678778
$ws.Cells.Item(25, 19).Value = 0.047323915106912
$ws.Cells.Item(25, 20).Value = 0.03574981517167066
# Row 26
$ws.Cells.Item(26, 7).Value = 124.2078576666667
$ws.Cells.Item(26, 8).Value = 372.623573
$ws.Cells.Item(26, 9).Value = 0.2163479508088675
$ws.Cells.Item(26, 10).Value = 0.2272780738416368
$ws.Cells.Item(26, 13).Value = 0.5965113333333333
$ws.Cells.Item(26, 14).Value = 1.789534
$ws.Cells.Item(26, 15).Value = 0.09301499989760488
$ws.Cells.Item(26, 16).Value = 0.1003304197230327
$ws.Cells.Item(26, 17).Value = 74.09139478722021
$ws.Cells.Item(26, 18).Value = 666.8225530849819
$ws.Cells.Item(26, 19).Value = 0.02012360462233384
$ws.Cells.Item(26, 20).Value = 0.02280290454237385
